$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Collin Sexton -> Anthony Edwards
$ws.Range("A3").Value = "Anthony Edwards"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Minnesota Timberwolves"

# Row 7: Giannis Antetokounmpo -> Kyle Kuzma
$ws.Range("A7").Value = "Kyle Kuzma"
$ws.Range("B7").Value = "PF"
$ws.Range("C7").Value = "Washington Wizards"

# Row 14: Payton Pritchard -> Collin Sexton
$ws.Range("A14").Value = "Collin Sexton"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Utah Jazz"

# Row 15: Kyle Kuzma -> Harrison Barnes
$ws.Range("A15").Value = "Harrison Barnes"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "San Antonio Spurs"

# Row 16: Anthony Edwards -> Giannis Antetokounmpo
$ws.Range("A16").Value = "Giannis Antetokounmpo"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Milwaukee Bucks"
